$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 175.5
$ws.Range("I9").Value = 175.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 175.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -6.5
$ws.Range("N9").ClearContents()
$ws.Range("H86").Value = 2009.8889
$ws.Range("I86").Value = 2025.8334
$ws.Range("J86").Value = 1978
$ws.Range("K86").Value = 2025.8334
$ws.Range("L86").Value = 1978
$ws.Range("M86").Value = -902.8334
$ws.Range("N86").Value = -4224
$ws.Range("H89").Value = 2009.8889
$ws.Range("I89").Value = 2025.8334
$ws.Range("J89").Value = 1978
$ws.Range("K89").Value = 10129.167
$ws.Range("L89").Value = 9890
$ws.Range("M89").Value = -4513.166999999999
$ws.Range("N89").Value = -21122
$ws.Range("H138").Value = 2566.7856
$ws.Range("I138").Value = 2366.3809
$ws.Range("J138").Value = 2767.1904
$ws.Range("K138").Value = 7099.1427
$ws.Range("L138").Value = 8301.5712
$ws.Range("M138").Value = -1959.1427
$ws.Range("N138").Value = -18581.5712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19876.033
$ws.Range("I32").Value = 22255.959
$ws.Range("K32").Value = 22255.959
$ws.Range("M32").Value = -21968.959
$ws.Range("H45").Value = 3661.0938
$ws.Range("I45").Value = 1540.4286
$ws.Range("K45").Value = 1540.4286
$ws.Range("M45").Value = -1163.4286
$ws.Range("H110").Value = 1466.6
$ws.Range("I110").Value = 1502.6666
$ws.Range("K110").Value = 1502.6666
$ws.Range("M110").Value = 542.3334
$ws.Range("H132").Value = 41510.883
$ws.Range("I132").Value = 53129.35
$ws.Range("K132").Value = 159388.05
$ws.Range("M132").Value = -156858.05

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2225.238
$ws.Range("I86").Value = 2179.6365
$ws.Range("J86").Value = 2275.4
$ws.Range("K86").Value = 2179.6365
$ws.Range("L86").Value = 2275.4
$ws.Range("M86").Value = -1056.6365
$ws.Range("N86").Value = -4521.4
$ws.Range("H89").Value = 2225.238
$ws.Range("I89").Value = 2179.6365
$ws.Range("J89").Value = 2275.4
$ws.Range("K89").Value = 10898.1825
$ws.Range("L89").Value = 11377
$ws.Range("M89").Value = -5282.182500000001
$ws.Range("N89").Value = -22609
$ws.Range("H105").Value = 120652.47
$ws.Range("I105").Value = 3051.6667
$ws.Range("J105").Value = 402894.4
$ws.Range("K105").Value = 3051.6667
$ws.Range("L105").Value = 402894.4
$ws.Range("M105").Value = -1304.6667
$ws.Range("N105").Value = -406388.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1375.5294
$ws.Range("I22").Value = 937.8333
$ws.Range("J22").Value = 1614.2727
$ws.Range("K22").Value = 937.8333
$ws.Range("L22").Value = 1614.2727
$ws.Range("M22").Value = -587.8333
$ws.Range("N22").Value = -2314.2727
$ws.Range("H31").Value = 4137.3076
$ws.Range("J31").Value = 7399.4
$ws.Range("L31").Value = 7399.4
$ws.Range("N31").Value = -7989.4
$ws.Range("H34").Value = 4137.3076
$ws.Range("J34").Value = 7399.4
$ws.Range("L34").Value = 7399.4
$ws.Range("N34").Value = -7803.4
$ws.Range("H99").Value = 4102
$ws.Range("J99").Value = 4828.3335
$ws.Range("L99").Value = 4828.3335
$ws.Range("N99").Value = -7824.3335
$ws.Range("H107").Value = 1727.7858
$ws.Range("I107").Value = 521
$ws.Range("K107").Value = 521
$ws.Range("M107").Value = 1399
$ws.Range("H126").Value = 4102
$ws.Range("J126").Value = 4828.3335
$ws.Range("L126").Value = 14485.0005
$ws.Range("N126").Value = -19425.0005
$ws.Range("H132").Value = 2880.818
$ws.Range("J132").Value = 2918.0476
$ws.Range("L132").Value = 8754.1428
$ws.Range("N132").Value = -13814.1428
$ws.Range("H134").Value = 168949.33
$ws.Range("I134").Value = 251724
$ws.Range("J134").Value = 3400
$ws.Range("K134").Value = 755172
$ws.Range("L134").Value = 10200
$ws.Range("M134").Value = -752637
$ws.Range("N134").Value = -15270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4033.6843
$ws.Range("I102").Value = 4645.909
$ws.Range("K102").Value = 4645.909
$ws.Range("M102").Value = -3023.909
$ws.Range("H122").Value = 2205.561
$ws.Range("I122").Value = 1632.2667
$ws.Range("K122").Value = 4896.800099999999
$ws.Range("M122").Value = -2446.800099999999
$ws.Range("H126").Value = 6346.1
$ws.Range("I126").Value = 5789.2856
$ws.Range("J126").Value = 7645.3335
$ws.Range("K126").Value = 17367.8568
$ws.Range("L126").Value = 22936.0005
$ws.Range("M126").Value = -14897.8568
$ws.Range("N126").Value = -27876.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8995.666999999999
$ws.Range("I16").Value = 8997.5
$ws.Range("J16").Value = 8992
$ws.Range("K16").Value = 8997.5
$ws.Range("L16").Value = 8992
$ws.Range("M16").Value = -8827.5
$ws.Range("N16").Value = -9332
$ws.Range("H46").Value = 46259.6
$ws.Range("I46").Value = 74499.664
$ws.Range("K46").Value = 74499.664
$ws.Range("M46").Value = -74311.664
$ws.Range("H55").Value = 414.6875
$ws.Range("I55").Value = 148.75
$ws.Range("J55").Value = 1212.5
$ws.Range("K55").Value = 148.75
$ws.Range("L55").Value = 1212.5
$ws.Range("M55").Value = 24.25
$ws.Range("N55").Value = -1558.5
$ws.Range("H61").Value = 1930.0526
$ws.Range("I61").Value = 1377.7333
$ws.Range("K61").Value = 1377.7333
$ws.Range("M61").Value = -1175.7333
$ws.Range("H82").Value = 2479.5264
$ws.Range("I82").Value = 1018.1111
$ws.Range("K82").Value = 1018.1111
$ws.Range("M82").Value = -657.1111
$ws.Range("H85").Value = 2479.5264
$ws.Range("I85").Value = 1018.1111
$ws.Range("K85").Value = 1018.1111
$ws.Range("M85").Value = 229.8889
$ws.Range("H93").Value = 5416.6665
$ws.Range("I93").Value = 5000
$ws.Range("K93").Value = 5000
$ws.Range("M93").Value = -3752
$ws.Range("H113").Value = 1930.0526
$ws.Range("I113").Value = 1377.7333
$ws.Range("K113").Value = 1377.7333
$ws.Range("M113").Value = 792.2666999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 79999
$ws.Range("J27").Value = 79999
$ws.Range("L27").Value = 79999
$ws.Range("N27").Value = -80137
